$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36, pushing existing rows 36-127 down to 37-128.
$ws.Rows("36:36").Insert()

# Populate the newly inserted row 36 with the new weekly record.
$ws.Range("A36").Value = 1
$ws.Range("B36").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C36").Value = "Arica y Parinacota"
$ws.Range("D36").Value = 44708
$ws.Range("D36").NumberFormat = $ws.Range("D37").NumberFormat
$ws.Range("E36").Value = 15
$ws.Range("F36").Value = "Fruta"
$ws.Range("G36").Value = 100108
$ws.Range("H36").Value = "Tropicales y subtropicales"
$ws.Range("I36").Value = 100108003
$ws.Range("J36").Value = "Maracuyá"
$ws.Range("K36").Value = "Sin especificar"
$ws.Range("L36").Value = "Primera"
$ws.Range("M36").Value = 140
$ws.Range("N36").Value = 22000
$ws.Range("O36").Value = 23000
$ws.Range("P36").Value = 22500
$ws.Range("Q36").Value = "$/caja 20 kilos"
$ws.Range("R36").Value = "Región de Arica y Parinacota"
$ws.Range("S36").Value = 1125
$ws.Range("T36").Value = 20
